# MazeLayout.xlsx - "Route finding partially working: tests need some revision"
#
# Adds a helper index column (V) alongside the existing maze grid (A:U,
# rows 1-20), and a helper index row (21) underneath it - both 0-based
# counters used by the route-finding logic under test - then moves the
# active selection to C11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column V (22): rows 1-20 hold a 0-based row index (0..19)
for ($r = 1; $r -le 20; $r++) {
    $ws.Cells.Item($r, 22).Value = $r - 1
}

# New row 21: columns A-U (1-21) hold a 0-based column index (0..20)
for ($c = 1; $c -le 21; $c++) {
    $ws.Cells.Item(21, $c).Value = $c - 1
}

# Narrow the new helper column to fit its single/double-digit contents
$ws.Columns.Item(22).ColumnWidth = 3 - 0.8333333333333334

# Move the active selection
$null = $ws.Range("C11").Select()
